$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row heights (rows 1-3 get new/changed explicit custom heights) ---
$ws.Rows.Item(1).RowHeight = 66
$ws.Rows.Item(2).RowHeight = 14.25
$ws.Rows.Item(3).RowHeight = 14.25

# --- New column N: copy formatting from column M (same row) so the new
#     cells inherit the matching number format / font / border / alignment,
#     then overwrite with the 2023 data values. ---

# Row 3 (blank separator row, just formatting, no value)
$ws.Range("M3").Copy($ws.Range("N3"))

# Row 4 (year header row)
$ws.Range("M4").Copy($ws.Range("N4"))
$ws.Range("N4").Value = 2023

# Row 5
$ws.Range("M5").Copy($ws.Range("N5"))
$ws.Range("N5").Value = 0

# Row 6
$ws.Range("M6").Copy($ws.Range("N6"))
$ws.Range("N6").Value = 48.5

# Row 7 - also bump L7's own format (picks up the one-decimal numeric format
# used elsewhere in the table, same as column M/K of the row below it)
$ws.Range("M8").Copy($ws.Range("L7"))
$ws.Range("L7").Value = 11.673077354810609

$ws.Range("M7").Copy($ws.Range("N7"))
$ws.Range("N7").Value = 23.2

# Row 8
$ws.Range("M8").Copy($ws.Range("N8"))
$ws.Range("N8").Value = 19.3

# Row 9
$ws.Range("M9").Copy($ws.Range("N9"))
$ws.Range("N9").Value = 9.1

$excel.CutCopyMode = $false

# --- Selection: clear the old "N7" selection left over from before the
#     column was added, back to the top-left cell. ---
$ws.Range("A1").Select()
